$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) The "_GoBack" bookmark currently sits between "...pick up changes "
#    and "then we need to update the copy in ...". It needs to move to
#    the very end of the new "ordering notation" paragraph we add below.
#    Delete it here; it gets re-created later in its new location.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2) Merge the run boundary left behind by the bookmark so the text
#    "...pick up changes then we need to update the copy in " is a
#    single contiguous run (matches the target canonical XML).
# ------------------------------------------------------------------
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute(
    "pick up changes then we need to update the copy in ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "pick up changes then we need to update the copy in ",
    2) | Out-Null

# ------------------------------------------------------------------
# 3) Append the new changelog entries (12/16/13) after the existing
#    "Added uri encoding..." paragraph, right before the trailing
#    empty paragraph that precedes the section break. The new, final
#    paragraph carries the relocated "_GoBack" bookmark at its end.
# ------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$insertionPoint = $lastPara.Range

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$newBody = "<w:p $wNs/>"
$newBody += "<w:p $wNs><w:r><w:t>12/16/13</w:t></w:r></w:p>"
$newBody += "<w:p $wNs>"
$newBody += "<w:r><w:t>Added the ability to order the questions by using the notation &lt;&lt;n&gt;&gt; at the start of each question. Have to order all the questions or none. If there are any errors in the notation will not attempt to order the questions.</w:t></w:r>"
$newBody += "<w:bookmarkStart w:id=`"1000`" w:name=`"_GoBack`"/>"
$newBody += "<w:bookmarkEnd w:id=`"1000`"/>"
$newBody += "</w:p>"

$packageXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' +
    $newBody +
    '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$insertionPoint.InsertXML($packageXml)

Write-Output "Edit applied."
